$wb = $excel.ActiveWorkbook

# Move the "总计" sheet so it becomes the first sheet (before "2022-Q2").
$summarySheet = $wb.Worksheets.Item("总计")
$firstSheet = $wb.Worksheets.Item(1)
$summarySheet.Move($firstSheet)
